# Add "PMID" column to the "studies" sheet and "notes" column to the
# "counts" sheet, matching the commit "added PMID and notes columns to
# data structure".

$wb = $excel.ActiveWorkbook

$studies = $wb.Worksheets.Item("studies")
$counts  = $wb.Worksheets.Item("counts")

# studies: header row 1 currently spans A1:G1 -> add H1 = "PMID"
# Match formatting of the adjoining header cell (G1) so the new header
# cell picks up the same font (black Aptos Narrow) used by the rest of
# the header row.
$studies.Range("H1").Value = "PMID"
$studies.Range("H1").Font.Color = $studies.Range("G1").Font.Color

# counts: header row 1 currently spans A1:E1 -> add F1 = "notes"
$counts.Range("F1").Value = "notes"

# Move selection / active cell to reflect the newly added column on each
# sheet, then finish with "counts" as the active sheet/tab.
$studies.Range("H2").Select()
$counts.Range("F2").Select()
$counts.Activate()
